$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.582
$ws.Range("D4").Value = -7.759
$ws.Range("E4").Value = 13.265

$ws.Range("D5").Value = -8.206999999999999

$ws.Range("A7").Value = -20.987

$ws.Range("D8").Value = -7.896000000000001

$ws.Range("E9").Value = 12.946

$ws.Range("A16").Value = -20.654
$ws.Range("D16").Value = -8.463999999999999

$ws.Range("E18").Value = 13.19
